$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column from 2023-09-15 (45184) to 2023-09-16 (45185)
# for rows 2 through 5 (serial date values, matching the existing numeric cell values).
$ws.Range("C2:C5").Value = 45185
